# Auto-generated Excel COM-interop script
# Applies cached-value refresh to H:N columns (currentAveragePrice* / LevePrice* / LeveProfit*)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, per the scheduled market-data runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

$ws.Range("H40").Value = 3750.625
$ws.Range("I40").Value = 3075.5715
$ws.Range("J40").Value = 5039.364
$ws.Range("K40").Value = 3075.5715
$ws.Range("L40").Value = 5039.364
$ws.Range("M40").Value = -2900.5715
$ws.Range("N40").Value = -5389.364

$ws.Range("H41").Value = 1531.5385
$ws.Range("I41").Value = 915.5
$ws.Range("J41").Value = 2059.5715
$ws.Range("K41").Value = 915.5
$ws.Range("L41").Value = 2059.5715
$ws.Range("M41").Value = -475.5
$ws.Range("N41").Value = -2939.5715

$ws.Range("H87").Value = 95176.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 95176.5
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 95176.5
$ws.Range("N87").Value = -97672.5

$ws.Range("H90").Value = 95176.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 95176.5
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 285529.5
$ws.Range("N90").Value = -298009.5

$ws.Range("H100").Value = 4459.125
$ws.Range("I100").Value = 3953.2856
$ws.Range("J100").Value = 8000
$ws.Range("K100").Value = 3953.2856
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = -3412.2856
$ws.Range("N100").Value = -9082

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1254

$ws.Range("H132").Value = 13908.5625
$ws.Range("I132").Value = 13908.5625
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 41725.6875
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -39195.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7386.75
$ws.Range("I32").Value = 6345.8667
$ws.Range("J32").Value = 23000
$ws.Range("K32").Value = 6345.8667
$ws.Range("L32").Value = 23000
$ws.Range("M32").Value = -6058.8667
$ws.Range("N32").Value = -23574

$ws.Range("H58").Value = 99995
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 99995
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 99995
$ws.Range("N58").Value = -100855

$ws.Range("H61").Value = 1663.3684
$ws.Range("I61").Value = 1663.3684
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1663.3684
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1451.3684

$ws.Range("H74").Value = 5244.3125
$ws.Range("I74").Value = 5196.92
$ws.Range("J74").Value = 5413.5713
$ws.Range("K74").Value = 5196.92
$ws.Range("L74").Value = 5413.5713
$ws.Range("M74").Value = -4322.92
$ws.Range("N74").Value = -7161.5713

$ws.Range("H77").Value = 5244.3125
$ws.Range("I77").Value = 5196.92
$ws.Range("J77").Value = 5413.5713
$ws.Range("K77").Value = 25984.6
$ws.Range("L77").Value = 27067.8565
$ws.Range("M77").Value = -21616.6
$ws.Range("N77").Value = -35803.85649999999

$ws.Range("H133").Value = 36500
$ws.Range("I133").Value = 23000
$ws.Range("J133").Value = 50000
$ws.Range("K133").Value = 23000
$ws.Range("L133").Value = 50000
$ws.Range("M133").Value = -20470
$ws.Range("N133").Value = -55060

$ws.Range("H136").Value = 1663.3684
$ws.Range("I136").Value = 1663.3684
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4990.1052
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2440.1052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1400.7368
$ws.Range("I20").Value = 1428.5555
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 1428.5555
$ws.Range("L20").Value = 900
$ws.Range("M20").Value = -1181.5555

$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 498
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 2010
$ws.Range("I105").Value = 2010
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2010
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1839.2727
$ws.Range("I22").Value = 711.7143
$ws.Range("J22").Value = 3812.5
$ws.Range("K22").Value = 711.7143
$ws.Range("L22").Value = 3812.5
$ws.Range("M22").Value = -361.7143
$ws.Range("N22").Value = -4512.5

$ws.Range("H31").Value = 6898.72
$ws.Range("I31").Value = 3072.7
$ws.Range("J31").Value = 9449.4
$ws.Range("K31").Value = 3072.7
$ws.Range("L31").Value = 9449.4
$ws.Range("M31").Value = -2777.7

$ws.Range("H34").Value = 6898.72
$ws.Range("I34").Value = 3072.7
$ws.Range("J34").Value = 9449.4
$ws.Range("K34").Value = 3072.7
$ws.Range("L34").Value = 9449.4
$ws.Range("M34").Value = -2870.7

$ws.Range("H58").Value = 4144.25
$ws.Range("I58").Value = 2243.25
$ws.Range("J58").Value = 6045.25
$ws.Range("K58").Value = 2243.25
$ws.Range("L58").Value = 6045.25
$ws.Range("M58").Value = -2040.25
$ws.Range("N58").Value = -6451.25

$ws.Range("H107").Value = 1816.6
$ws.Range("I107").Value = 1770.75
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1770.75
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 149.25
$ws.Range("N107").Value = -5840

$ws.Range("H132").Value = 3792.4666
$ws.Range("I132").Value = 3997.125
$ws.Range("J132").Value = 3558.5715
$ws.Range("K132").Value = 11991.375
$ws.Range("L132").Value = 10675.7145
$ws.Range("M132").Value = -9461.375
$ws.Range("N132").Value = -15735.7145

$ws.Range("H134").Value = 3199.8
$ws.Range("I134").Value = 3574.75
$ws.Range("J134").Value = 1700
$ws.Range("K134").Value = 10724.25
$ws.Range("L134").Value = 5100
$ws.Range("M134").Value = -8189.25

$ws.Range("H136").Value = 4144.25
$ws.Range("I136").Value = 2243.25
$ws.Range("J136").Value = 6045.25
$ws.Range("K136").Value = 6729.75
$ws.Range("L136").Value = 18135.75
$ws.Range("M136").Value = -4179.75
$ws.Range("N136").Value = -23235.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 119998.336
$ws.Range("I37").Value = 120000
$ws.Range("J37").Value = 119998
$ws.Range("K37").Value = 360000
$ws.Range("L37").Value = 359994
$ws.Range("M37").Value = -359888
$ws.Range("N37").Value = -360218

$ws.Range("H122").Value = 793.3333
$ws.Range("I122").Value = 690
$ws.Range("J122").Value = 896.6667
$ws.Range("K122").Value = 6210
$ws.Range("L122").Value = 8070.0003
$ws.Range("M122").Value = -3760
$ws.Range("N122").Value = -12970.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5879.375
$ws.Range("I7").Value = 4507.3335
$ws.Range("J7").Value = 9995.5
$ws.Range("K7").Value = 4507.3335
$ws.Range("L7").Value = 9995.5
$ws.Range("M7").Value = -4395.3335

$ws.Range("H16").Value = 994
$ws.Range("I16").Value = 994
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 994
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -824

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 3995
$ws.Range("I122").Value = 3995
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11985
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9535

$ws.Range("H126").Value = 5879.375
$ws.Range("I126").Value = 4507.3335
$ws.Range("J126").Value = 9995.5
$ws.Range("K126").Value = 13522.0005
$ws.Range("L126").Value = 29986.5
$ws.Range("M126").Value = -11052.0005

$ws.Range("H132").Value = 6328.875
$ws.Range("I132").Value = 2525.2
$ws.Range("J132").Value = 12668.333
$ws.Range("K132").Value = 7575.599999999999
$ws.Range("L132").Value = 38004.999
$ws.Range("M132").Value = -5045.599999999999

$ws.Range("H136").Value = 2462.4285
$ws.Range("I136").Value = 2317.3333
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 6951.999899999999
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -4401.999899999999
$ws.Range("N136").Value = -15099

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 816.0714
$ws.Range("I2").Value = 872.6923
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 872.6923
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = -760.6923

$ws.Range("H81").Value = 994.5
$ws.Range("I81").Value = 994.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1989
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -928

$ws.Range("H84").Value = 994.5
$ws.Range("I84").Value = 994.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9945
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4641

$ws.Range("H132").Value = 1518.2632
$ws.Range("I132").Value = 1160.5714
$ws.Range("J132").Value = 2519.8
$ws.Range("K132").Value = 3481.7142
$ws.Range("L132").Value = 7559.400000000001
$ws.Range("M132").Value = -951.7142000000003
$ws.Range("N132").Value = -12619.4

$ws.Range("H136").Value = 3368.3333
$ws.Range("I136").Value = 999
$ws.Range("J136").Value = 4553
$ws.Range("K136").Value = 999
$ws.Range("L136").Value = 13659
$ws.Range("M136").Value = -447
